# 24.11.2024 - kilka poprawek, powtorne importy, dodawanie w jakosci itp

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Summary formulas on row 1 (extend ranges to cover new rows) ---
$ws.Range("F1").Formula = "=COUNTA(C2:C200)"
$ws.Range("G1").Formula = "=SUM(D2:D200)"
$ws.Range("H1").Formula = "=G1/F1"

# --- Row 14: status + new remark in column E ---
$ws.Range("D14").Value = 80
$ws.Range("E14").Value = "Odkomentować zapis, dodać wykluczenie dodanych już grup roboczych i lokalizacji"

# --- Row 15 ---
$ws.Range("D15").Value = 100

# --- Row 18: status + new remark in column E ---
$ws.Range("D18").Value = 100
$ws.Range("E18").Value = "Wprowadzića tą metodę do innych okien, prawdopodobnie wszystkich"

# --- Row 20 ---
$ws.Range("D20").Value = 100

# --- New data rows 25-37 ---
$ws.Range("B25").Value = "Raporty"
$ws.Range("C25").Value = "Dostosować szerokości kolumn w zestawieniach"
$ws.Range("D25").Value = 100

$ws.Range("B26").Value = "Pracownicy"
$ws.Range("C26").Value = "Możliwość powtórnego wprowadzenia danych (skasowanie poprzednich i import nowych) "
$ws.Range("D26").Value = 100

$ws.Range("B27").Value = "Błędy produkcji"
$ws.Range("C27").Value = "Możliwość powtórnego wprowadzenia danych (skasowanie poprzednich i import nowych) "
$ws.Range("D27").Value = 100

$ws.Range("B28").Value = "Nieobecności"
$ws.Range("C28").Value = "Możliwość powtórnego wprowadzenia danych (skasowanie poprzednich i import nowych) "
$ws.Range("D28").Value = 100

$ws.Range("B29").Value = "Direct"
$ws.Range("C29").Value = "Możliwość powtórnego wprowadzenia danych (skasowanie poprzednich i import nowych) "
$ws.Range("D29").Value = 100

$ws.Range("B30").Value = "Raportownie (szczeg)"
$ws.Range("C30").Value = "Możliwość powtórnego wprowadzenia danych (skasowanie poprzednich i import nowych) "
$ws.Range("D30").Value = 100

$ws.Range("B31").Value = "Raportownie (total)"
$ws.Range("C31").Value = "Możliwość powtórnego wprowadzenia danych (skasowanie poprzednich i import nowych) "
$ws.Range("D31").Value = 100

$ws.Range("B32").Value = "Jakość "
$ws.Range("C32").Value = "Możliwość powtórnego wprowadzenia danych (skasowanie poprzednich i import nowych) "
$ws.Range("D32").Value = 100

$ws.Range("B33").Value = "Korekta Indirect"
$ws.Range("C33").Value = "Możliwość powtórnego wprowadzenia danych (skasowanie poprzednich i import nowych) "
$ws.Range("D33").Value = 100

$ws.Range("B34").Value = "Raporty"
$ws.Range("C34").Value = "poprawa raportu z danymi dla pracowników. Formatowanie komórek, wyrównanie, przedstawienie danych"
$ws.Range("D34").Value = 100

$ws.Range("B35").Value = "KPI magazynu "
$ws.Range("C35").Value = "Dodać możliwość kasowania jednego rekordu"
$ws.Range("D35").Value = 0

$ws.Range("B36").Value = "Jakość "
$ws.Range("C36").Value = "Sprawdzić poprawność filtrowania danych w grupach roboczych. Są różnie zapisanete same dane np. ""3013+3015"" i ""3013 + 3015"""
$ws.Range("D36").Value = 0
$ws.Range("A36:D36").RowHeight = 15.75

$ws.Range("B37").Value = "Instalator"
$ws.Range("C37").Value = "Stworzyć wersję .exe oraz doprowadzić by obrazy wyświetlały się i były pobierane z właściwej lokalizacji. "
$ws.Range("D37").Value = 100

# --- Trailing index-only rows 38-43 ---
for ($r = 38; $r -le 43; $r++) {
    $ws.Cells.Item($r, 1).Value = $r - 1
}

Write-Output "Data rows written"
